# Applies the "Controller para Usuarios y Productos creados." edit:
#  1) After the "Crear las clases en el paquete entidades (BackEnd)" bullet
#     (numId=5), insert a new bullet "Crear los Servicios, Repository y las
#     Interfaces para cada clase. (BackEnd)" followed by a new empty bullet.
#  2) Rewrite the existing "Crear los Servicios, Repository y las Interfaces
#     para cada clase. (BackEnd)" bullet (numId=7, ilvl=1) so it reads
#     "Crear Controller para cada clase. (BackEnd)".

$d = $word.ActiveDocument
$wordNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14Ns  = "http://schemas.microsoft.com/office/word/2010/wordml"

# NOTE: order matters. The text "Crear los Servicios, Repository y las
# Interfaces para cada clase. (BackEnd)" exists today only once (in the
# numId=7 bullet near "Lista de Tareas pendientes"), but step 2 below
# would *introduce* a second, earlier copy of that same sentence (in the
# numId=5 list). So we rewrite the numId=7 bullet FIRST, while the search
# text is still unique, and only afterwards insert the new numId=5
# bullets that duplicate the sentence.

# --- Step 1: locate the old "Crear los Servicios..." (numId=7) paragraph
#     and replace its content with "Crear Controller para cada clase."  ---

$targetRange = $d.Content
$null = $targetRange.Find.Execute("Crear los Servicios, Repository y las Interfaces para cada clase. (BackEnd)")
$targetRange.Expand(4) | Out-Null

$insertXml2 = @"
<w:p xmlns:w='$wordNs' xmlns:w14='$w14Ns' w14:paraId="09A90736" w14:textId="2DB8E2F0" w:rsidR="00EC683D" w:rsidRDefault="00624C8A" w:rsidP="00624C8A"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>Crear</w:t></w:r><w:r><w:t xml:space="preserve"> Controller</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00A82BC2"><w:t xml:space="preserve">para </w:t></w:r><w:r w:rsidR="0003455B"><w:t>cada</w:t></w:r><w:r w:rsidR="00A82BC2"><w:t xml:space="preserve"> clase.</w:t></w:r><w:r w:rsidR="00A82BC2" w:rsidRPr="00A82BC2"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00A82BC2"><w:t>(BackEnd)</w:t></w:r></w:p>
"@

$targetRange.InsertXML($insertXml2)

# --- Step 2: locate the "Crear las clases..." paragraph and append the
#     two new list paragraphs right after it. -----------------------------

$anchorRange = $d.Content
$null = $anchorRange.Find.Execute("Crear las clases en el paquete entidades (BackEnd)")
$anchorRange.Expand(4) | Out-Null

$insertXml1 = @"
<w:p xmlns:w='$wordNs' xmlns:w14='$w14Ns' w14:paraId="7D173761" w14:textId="2F2BF4E3" w:rsidR="00A82BC2" w:rsidRDefault="00A82BC2" w:rsidP="00D720C7"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Crear las clases en el paquete entidades (BackEnd)</w:t></w:r></w:p><w:p xmlns:w='$wordNs'><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Crear los Servicios, Repository y las Interfaces para cada clase.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(BackEnd)</w:t></w:r></w:p><w:p xmlns:w='$wordNs'><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr></w:p>
"@

$anchorRange.InsertXML($insertXml1)

Write-Output "Done."
